$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a zero-padded numeric-looking code ("002"); force text so Excel
# doesn't coerce it into the number 2, then drop the temporary number
# format back to Normal so no stray style sticks to the cell.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").Style = "Normal"

$ws.Range("N2").Value = "2020-06-30 00:00:00"

$ws.Range("O2").Value = -193769054.76
$ws.Range("P2").Value = 2522290245.31
$ws.Range("Q2").Value = 2707227729.75
$ws.Range("R2").Value = -56.9819255974
$ws.Range("S2").Value = 1697174134.51
$ws.Range("T2").Value = 1697174134.51
$ws.Range("U2").Value = -66.9044522617
$ws.Range("V2").Value = 714629910.51
$ws.Range("W2").Value = 127634961.88
$ws.Range("X2").Value = 116932073.19
$ws.Range("Y2").Value = -182709136.9
$ws.Range("Z2").Value = -170233935.58
$ws.Range("AA2").Value = 38800941.06
$ws.Range("AG2").Value = 50856649.66
$ws.Range("AP2").Value = -60.6057112233
$ws.Range("AQ2").Value = -260.503626672388
$ws.Range("AR2").Value = -385.62
$ws.Range("AS2").Value = -224704184.66
$ws.Range("AT2").Value = -479.925352862166
